# Apply scrims_actualizado.xlsx update: append new scrim-result rows
# to 5 worksheets (new matches logged 2025-07-25).
$wb = $excel.ActiveWorkbook

# Resolve every worksheet referenced below (as a direct target sheet
# or as a style-template donor for the cross-sheet "Equipo 2" copy).
$ws_Double_Swoosh = $wb.Worksheets.Item("Double Swoosh")
$ws_Dry_Season = $wb.Worksheets.Item("Dry Season")
$ws_Hot_Potato = $wb.Worksheets.Item("Hot Potato")
$ws_Layer_Cake = $wb.Worksheets.Item("Layer Cake")
$ws_New_Horizons = $wb.Worksheets.Item("New Horizons")

# --- Double Swoosh: A3:N9 -> A3:N13 (4 new rows) ---
$ws_Double_Swoosh.Range("A4:N4").Copy()
$ws_Double_Swoosh.Range("A10:N10").PasteSpecial(-4122)
$ws_Double_Swoosh.Range("A10").Value = "LOU"
$ws_Double_Swoosh.Range("B10").Value = "TARA"
$ws_Double_Swoosh.Range("C10").Value = "BULL"
$ws_Double_Swoosh.Range("D10").Value = "GRAY"
$ws_Double_Swoosh.Range("E10").Value = "JACKY"
$ws_Double_Swoosh.Range("F10").Value = "DRACO"
$ws_Double_Swoosh.Range("G10").Value = "Equipo 1"
$ws_Double_Swoosh.Range("H10").Value = "GEN|cookie"
$ws_Double_Swoosh.Range("I10").Value = "GEN|Moding"
$ws_Double_Swoosh.Range("J10").Value = "GEN|BONOX2"
$ws_Double_Swoosh.Range("K10").Value = "FZ|Mira"
$ws_Double_Swoosh.Range("L10").Value = "FZ|Danshari"
$ws_Double_Swoosh.Range("M10").Value = "FZ|Toridesu"
$ws_Double_Swoosh.Range("N10").Value = "20250725T132220.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_Double_Swoosh.Range("A11:N11").PasteSpecial(-4122)
$ws_Double_Swoosh.Range("A11").Value = "JACKY"
$ws_Double_Swoosh.Range("B11").Value = "TARA"
$ws_Double_Swoosh.Range("C11").Value = "GRAY"
$ws_Double_Swoosh.Range("D11").Value = "LILY"
$ws_Double_Swoosh.Range("E11").Value = "LUMI"
$ws_Double_Swoosh.Range("F11").Value = "ROSA"
$ws_Double_Swoosh.Range("G11").Value = "Equipo 2"
$ws_Double_Swoosh.Range("H11").Value = "GEN|Moding"
$ws_Double_Swoosh.Range("I11").Value = "GEN|BONOX2"
$ws_Double_Swoosh.Range("J11").Value = "GEN|cookie"
$ws_Double_Swoosh.Range("K11").Value = "FZ|Toridesu"
$ws_Double_Swoosh.Range("L11").Value = "FZ|Mira"
$ws_Double_Swoosh.Range("M11").Value = "FZ|Danshari"
$ws_Double_Swoosh.Range("N11").Value = "20250725T131650.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_Double_Swoosh.Range("A12:N12").PasteSpecial(-4122)
$ws_Double_Swoosh.Range("A12").Value = "JACKY"
$ws_Double_Swoosh.Range("B12").Value = "TARA"
$ws_Double_Swoosh.Range("C12").Value = "GRAY"
$ws_Double_Swoosh.Range("D12").Value = "LILY"
$ws_Double_Swoosh.Range("E12").Value = "LUMI"
$ws_Double_Swoosh.Range("F12").Value = "ROSA"
$ws_Double_Swoosh.Range("G12").Value = "Equipo 2"
$ws_Double_Swoosh.Range("H12").Value = "GEN|Moding"
$ws_Double_Swoosh.Range("I12").Value = "GEN|BONOX2"
$ws_Double_Swoosh.Range("J12").Value = "GEN|cookie"
$ws_Double_Swoosh.Range("K12").Value = "FZ|Toridesu"
$ws_Double_Swoosh.Range("L12").Value = "FZ|Mira"
$ws_Double_Swoosh.Range("M12").Value = "FZ|Danshari"
$ws_Double_Swoosh.Range("N12").Value = "20250725T131450.000Z"

$ws_Double_Swoosh.Range("A4:N4").Copy()
$ws_Double_Swoosh.Range("A13:N13").PasteSpecial(-4122)
$ws_Double_Swoosh.Range("A13").Value = "LOU"
$ws_Double_Swoosh.Range("B13").Value = "TARA"
$ws_Double_Swoosh.Range("C13").Value = "BULL"
$ws_Double_Swoosh.Range("D13").Value = "GRAY"
$ws_Double_Swoosh.Range("E13").Value = "JACKY"
$ws_Double_Swoosh.Range("F13").Value = "DRACO"
$ws_Double_Swoosh.Range("G13").Value = "Equipo 1"
$ws_Double_Swoosh.Range("H13").Value = "GEN|cookie"
$ws_Double_Swoosh.Range("I13").Value = "GEN|Moding"
$ws_Double_Swoosh.Range("J13").Value = "GEN|BONOX2"
$ws_Double_Swoosh.Range("K13").Value = "FZ|Mira"
$ws_Double_Swoosh.Range("L13").Value = "FZ|Danshari"
$ws_Double_Swoosh.Range("M13").Value = "FZ|Toridesu"
$ws_Double_Swoosh.Range("N13").Value = "20250725T132435.000Z"

# --- New Horizons: A3:N69 -> A3:N77 (8 new rows) ---
$ws_New_Horizons.Range("A4:N4").Copy()
$ws_New_Horizons.Range("A70:N70").PasteSpecial(-4122)
$ws_New_Horizons.Range("A70").Value = "CHARLIE"
$ws_New_Horizons.Range("B70").Value = "KAZE"
$ws_New_Horizons.Range("C70").Value = "MEEPLE"
$ws_New_Horizons.Range("D70").Value = "BROCK"
$ws_New_Horizons.Range("E70").Value = "BUSTER"
$ws_New_Horizons.Range("F70").Value = "ANGELO"
$ws_New_Horizons.Range("G70").Value = "Equipo 1"
$ws_New_Horizons.Range("H70").Value = "NAVI|Achapi"
$ws_New_Horizons.Range("I70").Value = "NAVI|Ryohei"
$ws_New_Horizons.Range("J70").Value = "NAVI|Kuru"
$ws_New_Horizons.Range("K70").Value = "あの頃のしずく👍"
$ws_New_Horizons.Range("L70").Value = "ZETA|Levi"
$ws_New_Horizons.Range("M70").Value = "あの頃のしてたんぽ👍"
$ws_New_Horizons.Range("N70").Value = "20250725T132324.000Z"

$ws_New_Horizons.Range("A4:N4").Copy()
$ws_New_Horizons.Range("A71:N71").PasteSpecial(-4122)
$ws_New_Horizons.Range("A71").Value = "CHARLIE"
$ws_New_Horizons.Range("B71").Value = "KAZE"
$ws_New_Horizons.Range("C71").Value = "MEEPLE"
$ws_New_Horizons.Range("D71").Value = "BROCK"
$ws_New_Horizons.Range("E71").Value = "BUSTER"
$ws_New_Horizons.Range("F71").Value = "ANGELO"
$ws_New_Horizons.Range("G71").Value = "Equipo 1"
$ws_New_Horizons.Range("H71").Value = "NAVI|Achapi"
$ws_New_Horizons.Range("I71").Value = "NAVI|Ryohei"
$ws_New_Horizons.Range("J71").Value = "NAVI|Kuru"
$ws_New_Horizons.Range("K71").Value = "あの頃のしずく👍"
$ws_New_Horizons.Range("L71").Value = "ZETA|Levi"
$ws_New_Horizons.Range("M71").Value = "あの頃のしてたんぽ👍"
$ws_New_Horizons.Range("N71").Value = "20250725T132108.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_New_Horizons.Range("A72:N72").PasteSpecial(-4122)
$ws_New_Horizons.Range("A72").Value = "CHARLIE"
$ws_New_Horizons.Range("B72").Value = "KAZE"
$ws_New_Horizons.Range("C72").Value = "MEEPLE"
$ws_New_Horizons.Range("D72").Value = "BROCK"
$ws_New_Horizons.Range("E72").Value = "BUSTER"
$ws_New_Horizons.Range("F72").Value = "ANGELO"
$ws_New_Horizons.Range("G72").Value = "Equipo 2"
$ws_New_Horizons.Range("H72").Value = "NAVI|Achapi"
$ws_New_Horizons.Range("I72").Value = "NAVI|Ryohei"
$ws_New_Horizons.Range("J72").Value = "NAVI|Kuru"
$ws_New_Horizons.Range("K72").Value = "あの頃のしずく👍"
$ws_New_Horizons.Range("L72").Value = "ZETA|Levi"
$ws_New_Horizons.Range("M72").Value = "あの頃のしてたんぽ👍"
$ws_New_Horizons.Range("N72").Value = "20250725T131838.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_New_Horizons.Range("A73:N73").PasteSpecial(-4122)
$ws_New_Horizons.Range("A73").Value = "OLLIE"
$ws_New_Horizons.Range("B73").Value = "LUMI"
$ws_New_Horizons.Range("C73").Value = "CROW"
$ws_New_Horizons.Range("D73").Value = "DRACO"
$ws_New_Horizons.Range("E73").Value = "BERRY"
$ws_New_Horizons.Range("F73").Value = "BROCK"
$ws_New_Horizons.Range("G73").Value = "Equipo 2"
$ws_New_Horizons.Range("H73").Value = "RC|Battoman"
$ws_New_Horizons.Range("I73").Value = "MM"
$ws_New_Horizons.Range("J73").Value = "RC|Shu"
$ws_New_Horizons.Range("K73").Value = "CR|Moya"
$ws_New_Horizons.Range("L73").Value = "CR|Milkreo"
$ws_New_Horizons.Range("M73").Value = "Tensai 천재"
$ws_New_Horizons.Range("N73").Value = "20250725T132300.000Z"

$ws_New_Horizons.Range("A4:N4").Copy()
$ws_New_Horizons.Range("A74:N74").PasteSpecial(-4122)
$ws_New_Horizons.Range("A74").Value = "OLLIE"
$ws_New_Horizons.Range("B74").Value = "LUMI"
$ws_New_Horizons.Range("C74").Value = "CROW"
$ws_New_Horizons.Range("D74").Value = "DRACO"
$ws_New_Horizons.Range("E74").Value = "BERRY"
$ws_New_Horizons.Range("F74").Value = "BROCK"
$ws_New_Horizons.Range("G74").Value = "Equipo 1"
$ws_New_Horizons.Range("H74").Value = "RC|Battoman"
$ws_New_Horizons.Range("I74").Value = "MM"
$ws_New_Horizons.Range("J74").Value = "RC|Shu"
$ws_New_Horizons.Range("K74").Value = "CR|Moya"
$ws_New_Horizons.Range("L74").Value = "CR|Milkreo"
$ws_New_Horizons.Range("M74").Value = "Tensai 천재"
$ws_New_Horizons.Range("N74").Value = "20250725T131911.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_New_Horizons.Range("A75:N75").PasteSpecial(-4122)
$ws_New_Horizons.Range("A75").Value = "OLLIE"
$ws_New_Horizons.Range("B75").Value = "LUMI"
$ws_New_Horizons.Range("C75").Value = "CROW"
$ws_New_Horizons.Range("D75").Value = "DRACO"
$ws_New_Horizons.Range("E75").Value = "BERRY"
$ws_New_Horizons.Range("F75").Value = "BROCK"
$ws_New_Horizons.Range("G75").Value = "Equipo 2"
$ws_New_Horizons.Range("H75").Value = "RC|Battoman"
$ws_New_Horizons.Range("I75").Value = "MM"
$ws_New_Horizons.Range("J75").Value = "RC|Shu"
$ws_New_Horizons.Range("K75").Value = "CR|Moya"
$ws_New_Horizons.Range("L75").Value = "CR|Milkreo"
$ws_New_Horizons.Range("M75").Value = "Tensai 천재"
$ws_New_Horizons.Range("N75").Value = "20250725T131701.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_New_Horizons.Range("A76:N76").PasteSpecial(-4122)
$ws_New_Horizons.Range("A76").Value = "CHARLIE"
$ws_New_Horizons.Range("B76").Value = "OLLIE"
$ws_New_Horizons.Range("C76").Value = "GRIFF"
$ws_New_Horizons.Range("D76").Value = "DRACO"
$ws_New_Horizons.Range("E76").Value = "CORDELIUS"
$ws_New_Horizons.Range("F76").Value = "KIT"
$ws_New_Horizons.Range("G76").Value = "Equipo 2"
$ws_New_Horizons.Range("H76").Value = "MM"
$ws_New_Horizons.Range("I76").Value = "RC|Battoman"
$ws_New_Horizons.Range("J76").Value = "RC|Shu"
$ws_New_Horizons.Range("K76").Value = "CR|Moya"
$ws_New_Horizons.Range("L76").Value = "CR|Milkreo"
$ws_New_Horizons.Range("M76").Value = "Tensai 천재"
$ws_New_Horizons.Range("N76").Value = "20250725T131000.000Z"

$ws_New_Horizons.Range("A8:N8").Copy()
$ws_New_Horizons.Range("A77:N77").PasteSpecial(-4122)
$ws_New_Horizons.Range("A77").Value = "CHARLIE"
$ws_New_Horizons.Range("B77").Value = "OLLIE"
$ws_New_Horizons.Range("C77").Value = "GRIFF"
$ws_New_Horizons.Range("D77").Value = "DRACO"
$ws_New_Horizons.Range("E77").Value = "CORDELIUS"
$ws_New_Horizons.Range("F77").Value = "KIT"
$ws_New_Horizons.Range("G77").Value = "Equipo 2"
$ws_New_Horizons.Range("H77").Value = "MM"
$ws_New_Horizons.Range("I77").Value = "RC|Battoman"
$ws_New_Horizons.Range("J77").Value = "RC|Shu"
$ws_New_Horizons.Range("K77").Value = "CR|Moya"
$ws_New_Horizons.Range("L77").Value = "CR|Milkreo"
$ws_New_Horizons.Range("M77").Value = "Tensai 천재"
$ws_New_Horizons.Range("N77").Value = "20250725T130752.000Z"

# --- Hot Potato: A3:N85 -> A3:N93 (8 new rows) ---
$ws_Hot_Potato.Range("A6:N6").Copy()
$ws_Hot_Potato.Range("A86:N86").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A86").Value = "FRANK"
$ws_Hot_Potato.Range("B86").Value = "CROW"
$ws_Hot_Potato.Range("C86").Value = "KIT"
$ws_Hot_Potato.Range("D86").Value = "DRACO"
$ws_Hot_Potato.Range("E86").Value = "CORDELIUS"
$ws_Hot_Potato.Range("F86").Value = "CARL"
$ws_Hot_Potato.Range("G86").Value = "Equipo 1"
$ws_Hot_Potato.Range("H86").Value = "NAVI|Ryohei"
$ws_Hot_Potato.Range("I86").Value = "NAVI|Kuru"
$ws_Hot_Potato.Range("J86").Value = "NAVI|Achapi"
$ws_Hot_Potato.Range("K86").Value = "あの頃のしてたんぽ👍"
$ws_Hot_Potato.Range("L86").Value = "あの頃のしずく👍"
$ws_Hot_Potato.Range("M86").Value = "ZETA|Levi"
$ws_Hot_Potato.Range("N86").Value = "20250725T131035.000Z"

$ws_Hot_Potato.Range("A6:N6").Copy()
$ws_Hot_Potato.Range("A87:N87").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A87").Value = "FRANK"
$ws_Hot_Potato.Range("B87").Value = "CROW"
$ws_Hot_Potato.Range("C87").Value = "KIT"
$ws_Hot_Potato.Range("D87").Value = "DRACO"
$ws_Hot_Potato.Range("E87").Value = "CORDELIUS"
$ws_Hot_Potato.Range("F87").Value = "CARL"
$ws_Hot_Potato.Range("G87").Value = "Equipo 1"
$ws_Hot_Potato.Range("H87").Value = "NAVI|Ryohei"
$ws_Hot_Potato.Range("I87").Value = "NAVI|Kuru"
$ws_Hot_Potato.Range("J87").Value = "NAVI|Achapi"
$ws_Hot_Potato.Range("K87").Value = "あの頃のしてたんぽ👍"
$ws_Hot_Potato.Range("L87").Value = "あの頃のしずく👍"
$ws_Hot_Potato.Range("M87").Value = "ZETA|Levi"
$ws_Hot_Potato.Range("N87").Value = "20250725T130933.000Z"

$ws_Hot_Potato.Range("A4:N4").Copy()
$ws_Hot_Potato.Range("A88:N88").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A88").Value = "CARL"
$ws_Hot_Potato.Range("B88").Value = "CHARLIE"
$ws_Hot_Potato.Range("C88").Value = "GRIFF"
$ws_Hot_Potato.Range("D88").Value = "BULL"
$ws_Hot_Potato.Range("E88").Value = "BARLEY"
$ws_Hot_Potato.Range("F88").Value = "CROW"
$ws_Hot_Potato.Range("G88").Value = "Equipo 2"
$ws_Hot_Potato.Range("H88").Value = "NAVI|Ryohei"
$ws_Hot_Potato.Range("I88").Value = "NAVI|Kuru"
$ws_Hot_Potato.Range("J88").Value = "NAVI|Achapi"
$ws_Hot_Potato.Range("K88").Value = "あの頃のしてたんぽ👍"
$ws_Hot_Potato.Range("L88").Value = "あの頃のしずく👍"
$ws_Hot_Potato.Range("M88").Value = "ZETA|Levi"
$ws_Hot_Potato.Range("N88").Value = "20250725T130418.000Z"

$ws_Hot_Potato.Range("A4:N4").Copy()
$ws_Hot_Potato.Range("A89:N89").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A89").Value = "CARL"
$ws_Hot_Potato.Range("B89").Value = "CHARLIE"
$ws_Hot_Potato.Range("C89").Value = "GRIFF"
$ws_Hot_Potato.Range("D89").Value = "BULL"
$ws_Hot_Potato.Range("E89").Value = "BARLEY"
$ws_Hot_Potato.Range("F89").Value = "CROW"
$ws_Hot_Potato.Range("G89").Value = "Equipo 2"
$ws_Hot_Potato.Range("H89").Value = "NAVI|Ryohei"
$ws_Hot_Potato.Range("I89").Value = "NAVI|Kuru"
$ws_Hot_Potato.Range("J89").Value = "NAVI|Achapi"
$ws_Hot_Potato.Range("K89").Value = "あの頃のしてたんぽ👍"
$ws_Hot_Potato.Range("L89").Value = "あの頃のしずく👍"
$ws_Hot_Potato.Range("M89").Value = "ZETA|Levi"
$ws_Hot_Potato.Range("N89").Value = "20250725T130225.000Z"

$ws_Hot_Potato.Range("A4:N4").Copy()
$ws_Hot_Potato.Range("A90:N90").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A90").Value = "CARL"
$ws_Hot_Potato.Range("B90").Value = "RICO"
$ws_Hot_Potato.Range("C90").Value = "CHARLIE"
$ws_Hot_Potato.Range("D90").Value = "BERRY"
$ws_Hot_Potato.Range("E90").Value = "EMZ"
$ws_Hot_Potato.Range("F90").Value = "KAZE"
$ws_Hot_Potato.Range("G90").Value = "Equipo 2"
$ws_Hot_Potato.Range("H90").Value = "MM"
$ws_Hot_Potato.Range("I90").Value = "RC|Shu"
$ws_Hot_Potato.Range("J90").Value = "RC|Battoman"
$ws_Hot_Potato.Range("K90").Value = "CR|Milkreo"
$ws_Hot_Potato.Range("L90").Value = "CR|Moya"
$ws_Hot_Potato.Range("M90").Value = "Tensai 천재"
$ws_Hot_Potato.Range("N90").Value = "20250725T130042.000Z"

$ws_Hot_Potato.Range("A4:N4").Copy()
$ws_Hot_Potato.Range("A91:N91").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A91").Value = "CARL"
$ws_Hot_Potato.Range("B91").Value = "RICO"
$ws_Hot_Potato.Range("C91").Value = "CHARLIE"
$ws_Hot_Potato.Range("D91").Value = "BERRY"
$ws_Hot_Potato.Range("E91").Value = "EMZ"
$ws_Hot_Potato.Range("F91").Value = "KAZE"
$ws_Hot_Potato.Range("G91").Value = "Equipo 2"
$ws_Hot_Potato.Range("H91").Value = "MM"
$ws_Hot_Potato.Range("I91").Value = "RC|Shu"
$ws_Hot_Potato.Range("J91").Value = "RC|Battoman"
$ws_Hot_Potato.Range("K91").Value = "CR|Milkreo"
$ws_Hot_Potato.Range("L91").Value = "CR|Moya"
$ws_Hot_Potato.Range("M91").Value = "Tensai 천재"
$ws_Hot_Potato.Range("N91").Value = "20250725T125911.000Z"

$ws_Hot_Potato.Range("A4:N4").Copy()
$ws_Hot_Potato.Range("A92:N92").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A92").Value = "R-T"
$ws_Hot_Potato.Range("B92").Value = "BERRY"
$ws_Hot_Potato.Range("C92").Value = "BUZZ"
$ws_Hot_Potato.Range("D92").Value = "CARL"
$ws_Hot_Potato.Range("E92").Value = "BARLEY"
$ws_Hot_Potato.Range("F92").Value = "KAZE"
$ws_Hot_Potato.Range("G92").Value = "Equipo 2"
$ws_Hot_Potato.Range("H92").Value = "MM"
$ws_Hot_Potato.Range("I92").Value = "RC|Battoman"
$ws_Hot_Potato.Range("J92").Value = "RC|Shu"
$ws_Hot_Potato.Range("K92").Value = "CR|Moya"
$ws_Hot_Potato.Range("L92").Value = "CR|Milkreo"
$ws_Hot_Potato.Range("M92").Value = "Tensai 천재"
$ws_Hot_Potato.Range("N92").Value = "20250725T125452.000Z"

$ws_Hot_Potato.Range("A6:N6").Copy()
$ws_Hot_Potato.Range("A93:N93").PasteSpecial(-4122)
$ws_Hot_Potato.Range("A93").Value = "R-T"
$ws_Hot_Potato.Range("B93").Value = "BERRY"
$ws_Hot_Potato.Range("C93").Value = "BUZZ"
$ws_Hot_Potato.Range("D93").Value = "CARL"
$ws_Hot_Potato.Range("E93").Value = "BARLEY"
$ws_Hot_Potato.Range("F93").Value = "KAZE"
$ws_Hot_Potato.Range("G93").Value = "Equipo 1"
$ws_Hot_Potato.Range("H93").Value = "MM"
$ws_Hot_Potato.Range("I93").Value = "RC|Battoman"
$ws_Hot_Potato.Range("J93").Value = "RC|Shu"
$ws_Hot_Potato.Range("K93").Value = "CR|Moya"
$ws_Hot_Potato.Range("L93").Value = "CR|Milkreo"
$ws_Hot_Potato.Range("M93").Value = "Tensai 천재"
$ws_Hot_Potato.Range("N93").Value = "20250725T125319.000Z"

# --- Layer Cake: A3:N91 -> A3:N93 (2 new rows) ---
$ws_Layer_Cake.Range("A9:N9").Copy()
$ws_Layer_Cake.Range("A92:N92").PasteSpecial(-4122)
$ws_Layer_Cake.Range("A92").Value = "BUSTER"
$ws_Layer_Cake.Range("B92").Value = "KIT"
$ws_Layer_Cake.Range("C92").Value = "LUMI"
$ws_Layer_Cake.Range("D92").Value = "JAE-YONG"
$ws_Layer_Cake.Range("E92").Value = "MEEPLE"
$ws_Layer_Cake.Range("F92").Value = "BARLEY"
$ws_Layer_Cake.Range("G92").Value = "Equipo 1"
$ws_Layer_Cake.Range("H92").Value = "NAVI|Ryohei"
$ws_Layer_Cake.Range("I92").Value = "NAVI|Achapi"
$ws_Layer_Cake.Range("J92").Value = "NAVI|Kuru"
$ws_Layer_Cake.Range("K92").Value = "あの頃のしてたんぽ👍"
$ws_Layer_Cake.Range("L92").Value = "ZETA|Levi"
$ws_Layer_Cake.Range("M92").Value = "あの頃のしずく👍"
$ws_Layer_Cake.Range("N92").Value = "20250725T125631.000Z"

$ws_Layer_Cake.Range("A9:N9").Copy()
$ws_Layer_Cake.Range("A93:N93").PasteSpecial(-4122)
$ws_Layer_Cake.Range("A93").Value = "BUSTER"
$ws_Layer_Cake.Range("B93").Value = "KIT"
$ws_Layer_Cake.Range("C93").Value = "LUMI"
$ws_Layer_Cake.Range("D93").Value = "JAE-YONG"
$ws_Layer_Cake.Range("E93").Value = "MEEPLE"
$ws_Layer_Cake.Range("F93").Value = "BARLEY"
$ws_Layer_Cake.Range("G93").Value = "Equipo 1"
$ws_Layer_Cake.Range("H93").Value = "NAVI|Ryohei"
$ws_Layer_Cake.Range("I93").Value = "NAVI|Achapi"
$ws_Layer_Cake.Range("J93").Value = "NAVI|Kuru"
$ws_Layer_Cake.Range("K93").Value = "あの頃のしてたんぽ👍"
$ws_Layer_Cake.Range("L93").Value = "ZETA|Levi"
$ws_Layer_Cake.Range("M93").Value = "あの頃のしずく👍"
$ws_Layer_Cake.Range("N93").Value = "20250725T125411.000Z"

# --- Dry Season: A3:N47 -> A3:N52 (5 new rows) ---
$ws_Dry_Season.Range("A6:N6").Copy()
$ws_Dry_Season.Range("A48:N48").PasteSpecial(-4122)
$ws_Dry_Season.Range("A48").Value = "GRAY"
$ws_Dry_Season.Range("B48").Value = "LUMI"
$ws_Dry_Season.Range("C48").Value = "KAZE"
$ws_Dry_Season.Range("D48").Value = "SQUEAK"
$ws_Dry_Season.Range("E48").Value = "BELLE"
$ws_Dry_Season.Range("F48").Value = "MEEPLE"
$ws_Dry_Season.Range("G48").Value = "Equipo 1"
$ws_Dry_Season.Range("H48").Value = "GEN|cookie"
$ws_Dry_Season.Range("I48").Value = "GEN|Moding"
$ws_Dry_Season.Range("J48").Value = "GEN|BONOX2"
$ws_Dry_Season.Range("K48").Value = "FZ|Danshari"
$ws_Dry_Season.Range("L48").Value = "FZ|Mira"
$ws_Dry_Season.Range("M48").Value = "FZ|Toridesu"
$ws_Dry_Season.Range("N48").Value = "20250725T130846.000Z"

$ws_Dry_Season.Range("A6:N6").Copy()
$ws_Dry_Season.Range("A49:N49").PasteSpecial(-4122)
$ws_Dry_Season.Range("A49").Value = "GRAY"
$ws_Dry_Season.Range("B49").Value = "LUMI"
$ws_Dry_Season.Range("C49").Value = "KAZE"
$ws_Dry_Season.Range("D49").Value = "SQUEAK"
$ws_Dry_Season.Range("E49").Value = "BELLE"
$ws_Dry_Season.Range("F49").Value = "MEEPLE"
$ws_Dry_Season.Range("G49").Value = "Equipo 1"
$ws_Dry_Season.Range("H49").Value = "GEN|cookie"
$ws_Dry_Season.Range("I49").Value = "GEN|Moding"
$ws_Dry_Season.Range("J49").Value = "GEN|BONOX2"
$ws_Dry_Season.Range("K49").Value = "FZ|Danshari"
$ws_Dry_Season.Range("L49").Value = "FZ|Mira"
$ws_Dry_Season.Range("M49").Value = "FZ|Toridesu"
$ws_Dry_Season.Range("N49").Value = "20250725T130627.000Z"

$ws_Dry_Season.Range("A4:N4").Copy()
$ws_Dry_Season.Range("A50:N50").PasteSpecial(-4122)
$ws_Dry_Season.Range("A50").Value = "GRAY"
$ws_Dry_Season.Range("B50").Value = "LUMI"
$ws_Dry_Season.Range("C50").Value = "KAZE"
$ws_Dry_Season.Range("D50").Value = "SQUEAK"
$ws_Dry_Season.Range("E50").Value = "BELLE"
$ws_Dry_Season.Range("F50").Value = "MEEPLE"
$ws_Dry_Season.Range("G50").Value = "Equipo 2"
$ws_Dry_Season.Range("H50").Value = "GEN|cookie"
$ws_Dry_Season.Range("I50").Value = "GEN|Moding"
$ws_Dry_Season.Range("J50").Value = "GEN|BONOX2"
$ws_Dry_Season.Range("K50").Value = "FZ|Danshari"
$ws_Dry_Season.Range("L50").Value = "FZ|Mira"
$ws_Dry_Season.Range("M50").Value = "FZ|Toridesu"
$ws_Dry_Season.Range("N50").Value = "20250725T130423.000Z"

$ws_Dry_Season.Range("A6:N6").Copy()
$ws_Dry_Season.Range("A51:N51").PasteSpecial(-4122)
$ws_Dry_Season.Range("A51").Value = "BELLE"
$ws_Dry_Season.Range("B51").Value = "GRIFF"
$ws_Dry_Season.Range("C51").Value = "KAZE"
$ws_Dry_Season.Range("D51").Value = "BYRON"
$ws_Dry_Season.Range("E51").Value = "LUMI"
$ws_Dry_Season.Range("F51").Value = "HANK"
$ws_Dry_Season.Range("G51").Value = "Equipo 1"
$ws_Dry_Season.Range("H51").Value = "GEN|cookie"
$ws_Dry_Season.Range("I51").Value = "GEN|Moding"
$ws_Dry_Season.Range("J51").Value = "GEN|BONOX2"
$ws_Dry_Season.Range("K51").Value = "FZ|Mira"
$ws_Dry_Season.Range("L51").Value = "FZ|Danshari"
$ws_Dry_Season.Range("M51").Value = "FZ|Toridesu"
$ws_Dry_Season.Range("N51").Value = "20250725T125806.000Z"

$ws_Dry_Season.Range("A4:N4").Copy()
$ws_Dry_Season.Range("A52:N52").PasteSpecial(-4122)
$ws_Dry_Season.Range("A52").Value = "BELLE"
$ws_Dry_Season.Range("B52").Value = "GRIFF"
$ws_Dry_Season.Range("C52").Value = "KAZE"
$ws_Dry_Season.Range("D52").Value = "BYRON"
$ws_Dry_Season.Range("E52").Value = "LUMI"
$ws_Dry_Season.Range("F52").Value = "HANK"
$ws_Dry_Season.Range("G52").Value = "Equipo 2"
$ws_Dry_Season.Range("H52").Value = "GEN|cookie"
$ws_Dry_Season.Range("I52").Value = "GEN|Moding"
$ws_Dry_Season.Range("J52").Value = "GEN|BONOX2"
$ws_Dry_Season.Range("K52").Value = "FZ|Mira"
$ws_Dry_Season.Range("L52").Value = "FZ|Danshari"
$ws_Dry_Season.Range("M52").Value = "FZ|Toridesu"
$ws_Dry_Season.Range("N52").Value = "20250725T125547.000Z"
